# Apply cryptocurrency price/volume updates scraped on Thu Aug  1 03:29:07 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.999.38"
$ws.Range("E2").Value = "  -3.14%  "

# Row 3
$ws.Range("D3").Value = "3.189.07"
$ws.Range("E3").Value = "  -3.11%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.50%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.20%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.613"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.67%  "

# Row 8
$ws.Range("E8").Value = "  -0.02%  "

# Row 9
$ws.Range("D9").Value = "3.189.00"
$ws.Range("E9").Value = "  -3.09%  "

# Row 10
$ws.Range("E10").Value = "  -3.18%  "

# Row 11
$ws.Range("E11").Value = "  -0.27%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.387"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.52%  "

# Row 13
$ws.Range("D13").Value = "3.740.91"
$ws.Range("E13").Value = "  -3.22%  "

# Row 14
$ws.Range("E14").Value = "  -1.82%  "

# Row 15
$ws.Range("D15").Value = "64.123.25"
$ws.Range("E15").Value = "  -3.02%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.46%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000160"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.15%  "

# Row 18
$ws.Range("D18").Value = "3.194.18"
$ws.Range("E18").Value = "  -2.67%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "416.51"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.03%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.19%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.86"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.92%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.99%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.09%  "

# Row 24
$ws.Range("E24").Value = "  -2.37%  "

# Row 25
$ws.Range("E25").Value = "  +3.25%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.493"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.16%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000110"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.98%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.44%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.31%  "

# Row 30
$ws.Range("E30").Value = "  -5.93%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.95"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.71%  "

# Row 32
$ws.Range("E32").Value = "  +0.02%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.93%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.00%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.14%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "156.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.24%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.18%  "

# Row 38
$ws.Range("D38").Value = "2.737.15"
$ws.Range("E38").Value = "  -1.58%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.78%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "25.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.08%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.31%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.718"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.00%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.93"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.15%  "

# Row 44
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0630"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.44%  "

# Row 45
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.05%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.16"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.54%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "298.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.70%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0263"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.98%  "

# Row 49
$ws.Range("E49").Value = "  -8.97%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0996"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.52%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.06%  "
